$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so Excel does not
# reinterpret numeric-looking strings (e.g. "1.000", "242.40") as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated price / volume figures scraped by the Action run.
$ws.Range("D2").Value = "29.912.65"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.876.52"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "0.7420"
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").Value = "242.40"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "0.07215"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "24.68"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").Value = "0.08412"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "0.7510"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "5.430"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "1.871.17"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D16").Value = "29.920.16"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "6.081"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "247.42"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "0.000007855"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "2.126.35"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "8.025"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "0.1560"
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("D26").Value = "9.262"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "164.91"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "2.038"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "1.504"
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("D31").Value = "4.603"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").Value = "1.531"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "4.272"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").Value = "0.05314"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "0.7541"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "2.692"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "0.01964"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "2.757"
$ws.Range("D41").Value = "0.4502"
$ws.Range("D42").Value = "1.109.76"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "6.050"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "72.58"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "0.8562"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "103.45"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.862"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.625"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "9.461"
$ws.Range("E50").Value = "  -3.27%  "
$ws.Range("D51").Value = "2.024.49"
$ws.Range("E51").Value = "  -3.19%  "

# Restore the default (no explicit) cell style on column D so the
# workbook keeps matching its original formatting/style indices.
$ws.Range("D2:D51").Style = "Normal"

